$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Country names (column A) and counts (column B) reflecting the new search query results
$data = @(
    @('Countries/Regions', 'Count'),
    @('PEOPLES R CHINA', 752),
    @('USA', 237),
    @('INDIA', 157),
    @('BELGIUM', 145),
    @('SOUTH KOREA', 119),
    @('ENGLAND', 117),
    @('ITALY', 92),
    @('NETHERLANDS', 88),
    @('CANADA', 87),
    @('AUSTRALIA', 80),
    @('TURKEY', 71),
    @('GERMANY', 64),
    @('BRAZIL', 54),
    @('JAPAN', 52),
    @('SPAIN', 47),
    @('SAUDI ARABIA', 45),
    @('PAKISTAN', 44),
    @('SINGAPORE', 42),
    @('WALES', 39),
    @('IRAN', 38),
    @('MALAYSIA', 36),
    @('EGYPT', 35),
    @('FRANCE', 35),
    @('POLAND', 32),
    @('IRAQ', 25),
    @('RUSSIA', 23),
    @('TURKIYE', 23),
    @('INDONESIA', 20),
    @('TAIWAN', 19),
    @('GREECE', 16),
    @('DENMARK', 14),
    @('ROMANIA', 14),
    @('SCOTLAND', 14),
    @('CZECH REPUBLIC', 13),
    @('NORWAY', 13),
    @('SWEDEN', 13),
    @('THAILAND', 11),
    @('PORTUGAL', 10),
    @('QATAR', 10),
    @('SERBIA', 9),
    @('SWITZERLAND', 9),
    @('LITHUANIA', 8),
    @('NEW ZEALAND', 8),
    @('VIETNAM', 8),
    @('ETHIOPIA', 7),
    @('ISRAEL', 6),
    @('SOUTH AFRICA', 6),
    @('U ARAB EMIRATES', 6),
    @('ARGENTINA', 5),
    @('BANGLADESH', 5),
    @('COLOMBIA', 5),
    @('LATVIA', 5),
    @('LIBYA', 5),
    @('ALGERIA', 4),
    @('AUSTRIA', 4),
    @('NIGERIA', 4),
    @('TUNISIA', 4),
    @('CHILE', 3),
    @('CROATIA', 3),
    @('CYPRUS', 3),
    @('ECUADOR', 3),
    @('FINLAND', 3),
    @('JORDAN', 3),
    @('KOSOVO', 3),
    @('MALTA', 3),
    @('MEXICO', 3),
    @('OMAN', 3),
    @('IRELAND', 2),
    @('KENYA', 2),
    @('PHILIPPINES', 2),
    @('SLOVENIA', 2),
    @('YEMEN', 2),
    @('ALBANIA', 1),
    @('AZERBAIJAN', 1),
    @('BAHRAIN', 1),
    @('BULGARIA', 1),
    @('ESTONIA', 1),
    @('HUNGARY', 1),
    @('ICELAND', 1),
    @('KAZAKHSTAN', 1),
    @('KUWAIT', 1),
    @('LUXEMBOURG', 1),
    @('MALAWI', 1),
    @('PALESTINE', 1),
    @('PERU', 1),
    @('SLOVAKIA', 1),
    @('SRI LANKA', 1),
    @('TANZANIA', 1),
    @('UGANDA', 1)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $i + 1
    $ws.Cells.Item($rowNum, 1).Value = $data[$i][0]
    $ws.Cells.Item($rowNum, 2).Value = $data[$i][1]
}

# Row 91 no longer exists in the updated data - clear it so the used range shrinks to B90
$ws.Range("A91:B91").ClearContents()

$wb.Save()
